$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.772.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0637"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.24"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.632.45"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.858.48"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0770"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.27"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.787.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.37"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.895"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.108.59"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.55"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.803"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0109"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.09"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.71"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.01%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.19%  "
